$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells keep their text (string) representation exactly,
# matching the source workbook where every data cell is an inline string.
$cellUpdates = @{
    'D2' = '319.72'
    'E2' = '5.75%'
    'D3' = '49.44'
    'E3' = '14.17%'
    'D4' = '5.260'
    'E4' = '3.38%'
    'D5' = '0.07936'
    'E5' = '3.05%'
    'D6' = '4.578'
    'E6' = '3.69%'
    'D7' = '1.388'
    'E7' = '33.98%'
    'D8' = '1.638'
    'E8' = '1.19%'
    'D9' = '0.1297'
    'E9' = '3.48%'
    'D10' = '0.1969'
    'E10' = '6.31%'
    'D11' = '0.09484'
    'E11' = '3.52%'
    'D12' = '0.04607'
    'E12' = '10.46%'
    'D13' = '0.1046'
    'E13' = '-0.24%'
    'E14' = '2.12%'
    'D15' = '0.04173'
    'E15' = '0.46%'
    'D16' = '0.005902'
    'E16' = '2.43%'
    'B17' = 'LEO'
    'C17' = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
    'D17' = '3.343'
    'E17' = '0.07%'
    'B18' = 'BTSEToken'
    'C18' = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
    'D18' = '2.436'
    'E18' = '3.43%'
    'B19' = 'BitpandaEcosystemToken'
    'C19' = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
    'D19' = '0.3457'
    'E19' = '3.08%'
    'B20' = 'MCDex'
    'C20' = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
    'D20' = '8.211'
    'E20' = '-4.85%'
    'B21' = 'ProBitToken'
    'C21' = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
    'D21' = '0.1390'
    'E21' = '1.66%'
    'B22' = 'ZBToken'
    'C22' = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
    'D22' = '0.3090'
    'E22' = '-3.31%'
    'B23' = 'BitKan'
    'C23' = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
    'D23' = '0.001313'
    'E23' = '2.20%'
    'B24' = 'HotbitToken'
    'C24' = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
    'D24' = '0.004258'
    'E24' = '-4.72%'
    'D25' = '0.0001348'
    'E25' = '-0.10%'
    'D26' = '0.0003534'
    'D38' = '0.02689'
    'E38' = '9.47%'
    'D39' = '0.05791'
    'E39' = '9.94%'
    'D40' = '0.01092'
    'E40' = '83.12%'
    'D41' = '0.007994'
    'E41' = '4.25%'
    'D42' = '0.1442'
    'E42' = '7.07%'
    'D43' = '0.007679'
    'E43' = '4.41%'
    'D44' = '0.008463'
    'E44' = '11.91%'
    'D45' = '0.3196'
    'E45' = '6.24%'
    'D46' = '0.00006615'
    'E46' = '-1.33%'
    'D47' = '0.00000000749'
    'E47' = '-0.08%'
    'D48' = '0.05493'
    'E48' = '22.62%'
    'D49' = '0.003995'
    'E49' = '-4.91%'
    'D50' = '0.00002097'
    'E50' = '-0.08%'
    'D51' = '0.0001997'
    'E51' = '-0.08%'
}

foreach ($ref in $cellUpdates.Keys) {
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $cellUpdates[$ref]
}
